$wb = $excel.ActiveWorkbook

# --- Summary sheet value updates ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 50
$wsSummary.Range("E4").Value = 0

# --- Repayment schedule sheet value updates ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("I7").Value = 0
$wsRepay.Range("K7").Value = 887.72
$wsRepay.Range("P7").Value = 887.72

# --- Selection / active-tab updates ---
# Select the target cell on "Repayment schedule" first (temporarily activates it)
[void]$wsRepay.Range("L20").Select()

# Finally select the target cell on "Summary" - this activates Summary as the
# last-active sheet, matching the final workbook state (activeTab=1, tabSelected
# moves from Transactions to Summary).
[void]$wsSummary.Range("I5").Select()
